# KENGINE_MONDELEZUS_V1 template update:
# - Merge "Secondary Location" KPI params into the "Primary Location" cell text
#   (comma separated) for Candy/Chocolate, Chewing Gum and Biscuits rows, and
#   drop the now-redundant standalone "Secondary Location" row (D4 becomes a
#   duplicate of D2).
# - Row heights grow to fit the now-longer, wrapped text.
# - Selection / column-width bookkeeping to mirror the refreshed template.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KPI")

# --- Candy/Chocolate: merge Primary + Secondary location into D2, and make
#     D4 (previously the standalone "Secondary Location" row) mirror D2.
$ws.Range("D2").Value = "Candy/Chocolate Primary Location,Candy/Chocolate Secondary Location"
$ws.Range("D2").Copy($ws.Range("D4"))

# --- Chewing Gum: merge Primary + Secondary location into D3.
$ws.Range("D3").Value = "Chewing Gum Primary Location ,Chewing Gum Secondary Location"

# --- Biscuits (Cookies & Crackers): merge Primary + Secondary location into D5.
$ws.Range("D5").Value = "Biscuits (Cookies & Crackers) Primary Location,Biscuits (Cookies & Crackers) Secondary Location"

# Row heights grow because of the longer wrapped text.
$ws.Rows.Item(2).RowHeight = 28.35
$ws.Rows.Item(4).RowHeight = 28.1
$ws.Rows.Item(5).RowHeight = 28.35

# Cursor moved down one row after the edit.
$ws.Range("D4").Select()

$wsVisible = $wb.Worksheets.Item("Visible")
$wsVisible.Range("B8").Select()

$wsSetSize = $wb.Worksheets.Item("set size")
$wsSetSize.Range("D21").Select()

$ws.Activate()
